$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "B7"  = 5.048300000000001
    "A9"  = -21.91760000000001
    "B12" = 5.433799999999996
    "D15" = -8.382799999999998
    "A18" = -22.38570000000002
    "A20" = -19.11559999999999
    "B26" = 4.012200000000002
    "A27" = -21.70399999999997
    "B27" = 4.772700000000004
    "B29" = 5.229799999999998
    "B37" = 8.785100000000007
    "B38" = 4.3847
    "D38" = -9.013799999999998
    "D44" = -7.194300000000004
    "B51" = 5.3752
    "D51" = -7.5108
    "B55" = 5.168999999999998
    "D57" = -8.2498
    "D63" = -7.9762
    "A69" = -21.70129999999999
    "B69" = 5.657199999999996
    "B70" = 5.738699999999999
    "D70" = -6.841999999999998
    "A76" = -19.38749999999999
    "A82" = -22.111
    "B83" = 5.910000000000002
    "D99" = -7.7972
    "B102" = 8.337700000000005
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
